$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 66522
$ws.Range("B2").Value = "Stephany Ribeiro"
$ws.Range("C2").Value = "Recursos Humanos"
$ws.Range("D2").Value = "Viagem de negócios"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 45078
$ws.Range("G2").Value = 7765.27

# Row 3
$ws.Range("A3").Value = 58575
$ws.Range("B3").Value = "Thomas Moraes"
$ws.Range("C3").Value = "TI"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45085
$ws.Range("G3").Value = 5453.5

# Row 4
$ws.Range("A4").Value = 68558
$ws.Range("B4").Value = "Felipe Duarte"
$ws.Range("C4").Value = "Marketing"
$ws.Range("D4").Value = "Viagem de negócios"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45082
$ws.Range("G4").Value = 5773.76

# Row 5
$ws.Range("A5").Value = 19854
$ws.Range("B5").Value = "Sr. Davi Lucca Teixeira"
$ws.Range("C5").Value = "P&D"
$ws.Range("D5").Value = "Consulta médica"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 45082
$ws.Range("G5").Value = 9058

# Row 6
$ws.Range("A6").Value = 33847
$ws.Range("B6").Value = "Bryan Silva"
$ws.Range("C6").Value = "Engenharia"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 45089
$ws.Range("G6").Value = 11391.62

# Row 7
$ws.Range("A7").Value = 25463
$ws.Range("B7").Value = "Alexia Pereira"
$ws.Range("C7").Value = "Vendas"
$ws.Range("D7").Value = "Consulta médica"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45104
$ws.Range("G7").Value = 10105.43

# Row 8
$ws.Range("A8").Value = 20507
$ws.Range("B8").Value = "Ana Carolina Fogaça"
$ws.Range("C8").Value = "TI"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 45096
$ws.Range("G8").Value = 10082.31

# Row 9
$ws.Range("A9").Value = 85368
$ws.Range("B9").Value = "Helena Silveira"
$ws.Range("C9").Value = "Atendimento ao Cliente"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 45097
$ws.Range("G9").Value = 6223.66

# Row 10
$ws.Range("A10").Value = 81519
$ws.Range("B10").Value = "Lorena Vieira"
$ws.Range("C10").Value = "Operações"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 45095
$ws.Range("G10").Value = 7700.17

# Row 11
$ws.Range("A11").Value = 71119
$ws.Range("B11").Value = "Dra. Alícia Nogueira"
$ws.Range("C11").Value = "TI"
$ws.Range("D11").Value = "Consulta médica"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45083
$ws.Range("G11").Value = 5650.22
